$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.251.59'
$ws.Range('E2').Value = '  -2.92%  '
$ws.Range('D3').Value = '3.296.47'
$ws.Range('E3').Value = '  -3.67%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.68'
$ws.Range('E5').Value = '  -3.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.73'
$ws.Range('E6').Value = '  -8.67%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.297.98'
$ws.Range('E8').Value = '  -3.62%  '
$ws.Range('E9').Value = '  -3.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.89'
$ws.Range('E10').Value = '  -2.40%  '
$ws.Range('E11').Value = '  -5.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.406'
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D13').Value = '3.861.35'
$ws.Range('E13').Value = '  -3.72%  '
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.49'
$ws.Range('E15').Value = '  -7.54%  '
$ws.Range('D16').Value = '3.298.55'
$ws.Range('E16').Value = '  -3.57%  '
$ws.Range('E17').Value = '  -4.99%  '
$ws.Range('D18').Value = '60.247.33'
$ws.Range('E18').Value = '  -2.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.63'
$ws.Range('E20').Value = '  -5.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.53'
$ws.Range('E21').Value = '  -5.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '373.32'
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.99'
$ws.Range('E24').Value = '  -5.18%  '
$ws.Range('E25').Value = '  -7.18%  '
$ws.Range('D26').Value = '3.435.48'
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('E27').Value = '  -9.90%  '
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.03'
$ws.Range('E30').Value = '  -8.30%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  -5.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.29'
$ws.Range('E33').Value = '  -7.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.55'
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.25'
$ws.Range('E35').Value = '  -5.41%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.03'
$ws.Range('E36').Value = '  -8.70%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '165.45'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  -4.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.63'
$ws.Range('E39').Value = '  -4.92%  '
$ws.Range('D40').Value = '3.325.72'
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0721'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '25.61'
$ws.Range('E42').Value = '  -17.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.76'
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.745'
$ws.Range('E44').Value = '  -4.44%  '
$ws.Range('E45').Value = '  -4.24%  '
$ws.Range('E46').Value = '  -7.17%  '
$ws.Range('E47').Value = '  -6.70%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').Value = '2.324.42'
$ws.Range('E49').Value = '  -9.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.38'
$ws.Range('E50').Value = '  -6.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.47'
$ws.Range('E51').Value = '  -8.31%  '
